$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.442.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.519.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.37%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.07%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  +0.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.521.96'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.45%  '

$ws.Range("E10").Value = '  +0.60%  '

$ws.Range("E11").Value = '  -1.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("E13").Value = '  -0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("E15").Value = '  +1.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.990.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.350.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.531.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '328.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '644.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("E28").Value = '  +5.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.655.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.46%  '

$ws.Range("E30").Value = '  +5.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.996'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.69%  '

$ws.Range("E34").Value = '  +2.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("E37").Value = '  +1.27%  '

$ws.Range("E38").Value = '  +2.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '153.72'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("E40").Value = '  +0.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.91'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("E43").Value = '  +2.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '163.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.12%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0298'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.617'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.74%  '

$ws.Range("E51").Value = '  +1.57%  '
